$d = $word.ActiveDocument
$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Change 1 ---------------------------------------------------------
# Insert a new paragraph (carrying the <w:lastRenderedPageBreak/> that
# used to sit on the "3) Deploy the WAR file..." paragraph) right before
# that paragraph, with the new "Or <<GIT location download>>/..." text,
# then re-emit the original paragraph without the page-break marker.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Deploy the WAR file into TOMCAT server*") {
        $r = $p.Range
        $xml = '<w:p xmlns:w="' + $w + '">' +
                 '<w:r>' +
                   '<w:lastRenderedPageBreak/>' +
                   '<w:t xml:space="preserve"> Or &lt;&lt;GIT location download&gt;&gt;/</w:t>' +
                 '</w:r>' +
                 '<w:r><w:t>\SpringDemoGradleProject</w:t></w:r>' +
                 '<w:r><w:t>/</w:t></w:r>' +
                 '<w:r><w:t>SpringDemoGradleProject-0.0.1-SNAPSHOT</w:t></w:r>' +
                 '<w:r><w:t>.war</w:t></w:r>' +
               '</w:p>' +
               '<w:p xmlns:w="' + $w + '" w:rsidR="00B61296" w:rsidRDefault="00B61296">' +
                 '<w:r><w:t xml:space="preserve">3) Deploy the WAR file into TOMCAT server </w:t></w:r>' +
                 '<w:proofErr w:type="gramStart"/>' +
                 '<w:r><w:t xml:space="preserve">manually </w:t></w:r>' +
                 '<w:r w:rsidR="00F00D7E"><w:t xml:space="preserve"> ,</w:t></w:r>' +
                 '<w:proofErr w:type="gramEnd"/>' +
                 '<w:r w:rsidR="00F00D7E"><w:t xml:space="preserve"> start the Tomcat server after copy the WAR file </w:t></w:r>' +
               '</w:p>'
        $r.InsertXML($xml)
        break
    }
}

# --- Change 2 ---------------------------------------------------------
# Add a <w:lastRenderedPageBreak/> marker at the start of the
# "2.5)  Test using ..." paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Test using above URL*") {
        $r = $p.Range
        $xml = '<w:p xmlns:w="' + $w + '" w:rsidR="00F00D7E" w:rsidRDefault="00F00D7E" w:rsidP="00F00D7E">' +
                 '<w:r>' +
                   '<w:lastRenderedPageBreak/>' +
                   '<w:t xml:space="preserve">2.5)  Test using </w:t>' +
                 '</w:r>' +
                 '<w:r w:rsidR="000537B1"><w:t xml:space="preserve">above </w:t></w:r>' +
                 '<w:proofErr w:type="gramStart"/>' +
                 '<w:r w:rsidR="000537B1"><w:t>URL(</w:t></w:r>' +
                 '<w:proofErr w:type="gramEnd"/>' +
                 '<w:r w:rsidR="000537B1"><w:t xml:space="preserve"> Update with new URL ) and test the same from SOAP UI as well.</w:t></w:r>' +
               '</w:p>'
        $r.InsertXML($xml)
        break
    }
}
